$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab: SCD0297 -> SCD0018
$ws.Name = "SCD0018"

# Update TC_ID column (B2:B6) from "DGS-312" to "SCD0018-020"
$ws.Range("B2:B6").Value = "SCD0018-020"

# Column B needs to widen to fit the longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6

# Move the selection to match the saved view (active cell B7)
$ws.Range("B7").Select()
